# Results from June 11, 2020 12:01 AM run
# Replaces the disparities data table (A2:I12) with the refreshed,
# alphabetically-sorted dataset (now A2:I16): four new states added
# (Alabama, Arkansas, California, Colorado) and all previously-reported
# rows refreshed with the latest scraped values.

function Set-DataCell($ws, $r, $c, $type, $val) {
    if ($type -eq "N") {
        $ws.Cells.Item($r, $c).Value = $val
    } elseif ($type -eq "T") {
        $ws.Cells.Item($r, $c).NumberFormat = "@"
        $ws.Cells.Item($r, $c).Value = $val
        $ws.Cells.Item($r, $c).Style = "Normal"
    } elseif ($type -eq "D2") {
        $ws.Cells.Item($r, $c).NumberFormat = "YYYY-MM-DD"
        $ws.Cells.Item($r, $c).Value = $val
    } elseif ($type -eq "D3") {
        $ws.Cells.Item($r, $c).NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Cells.Item($r, $c).Value = $val
    } elseif ($type -eq "E") {
        $ws.Cells.Item($r, $c).Value = ""
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old table body (content + formatting) before rewriting it in
# the new sorted order.
$ws.Range("A2:I16").Clear()

$data = @(
    ,@('Alabama', 'D2', 43992, 'N', 21626, 'N', 739, 'N', 9221, 'N', 333, 'N', 42.64, 'N', 45.06, 'Success!')
    ,@('Arkansas', 'D2', 43992, 'N', 593, 'N', 11, 'N', 16, 'N', 0, 'T', '2.7', 'T', '0.0', 'Success!')
    ,@('California', 'D2', 43991, 'T', '97336', 'T', '4600', 'T', '4713', 'T', '451', 'T', '4.8', 'T', '9.8', 'Success!')
    ,@('California - San Diego', 'D3', 43991.99836937636, 'N', 8729, 'N', 301, 'N', 260, 'N', 9, 'N', 2.98, 'N', 2.99, 'Success!')
    ,@('Colorado', 'D2', 43992, 'N', 28499, 'N', 1573, 'N', 1676, 'N', 106, 'N', 5.88, 'N', 6.74, 'Success!')
    ,@('Florida', 'T', '2020-06-10', 'N', 65779, 'N', 2801, 'N', 12198, 'N', 558, 'N', 18.54, 'N', 19.92, 'Success!')
    ,@('Georgia', 'D2', 43992, 'N', 53980, 'N', 2329, 'N', 16965, 'N', 1123, 'N', 31.43, 'N', 48.22, 'Success!')
    ,@('Massachusetts', 'D3', 43992, 'N', 104156, 'N', 7454, 'N', 9729, 'N', 620, 'N', 9.34, 'N', 8.32, 'Success!')
    ,@('Michigan', 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'An error occurred. ... UnboundLocalError("local variable ''date_published'' referenced before assignment")')
    ,@('Minnesota', 'D2', 43992, 'N', 28869, 'N', 1236, 'N', 6342, 'N', 78, 'N', 21.97, 'N', 6.31, 'Success!')
    ,@('North Carolina', 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'An error occurred. ... ValueError(''Unable to extract date from table header.'')')
    ,@('Texas -- Bexar County', 'D2', 43992, 'N', 1805, 'N', 54, 'N', 280, 'N', 17, 'N', 15.51, 'N', 31.48, 'Success!')
    ,@('Virginia', 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'E', $null, 'An error occurred. ... URLError(TimeoutError(10060, ''A connection attempt failed because the connected party did not properly respond after a period of time, or established connection failed because connected host has failed to respond'', None, 10060, None))')
    ,@('Washington, DC', 'D3', 43991, 'N', 9474, 'N', 495, 'N', 4331, 'N', 367, 'N', 45.71, 'N', 74.14, 'Success!')
    ,@('Wisconsin -- Milwaukee', 'D2', 43992, 'N', 9161, 'N', 306, 'N', 2597, 'N', 129, 'N', 28.35, 'N', 1.41, 'Success!')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    Set-DataCell $ws $r 2 $row[1] $row[2]
    Set-DataCell $ws $r 3 $row[3] $row[4]
    Set-DataCell $ws $r 4 $row[5] $row[6]
    Set-DataCell $ws $r 5 $row[7] $row[8]
    Set-DataCell $ws $r 6 $row[9] $row[10]
    Set-DataCell $ws $r 7 $row[11] $row[12]
    Set-DataCell $ws $r 8 $row[13] $row[14]
    $ws.Cells.Item($r, 9).Value = $row[15]
    $r = $r + 1
}
